# Update the "December 2021 - Present" entry on Sheet1 to reflect an end date.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A11").Value = "December 2021 - April 2022"

# Move the active selection to A12, matching the author's cursor position after editing.
$ws.Range("A12").Select()
